$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Monday 2018-07-02 timesheet entry
$ws.Range("A10").Value = 43283
$ws.Range("B10").Value = 0.39583333333333331
$ws.Range("C10").Value = 0.20833333333333334
$ws.Range("D10").Value = 0.38541666666666669
$ws.Range("E10").Value = 0.20833333333333334
$ws.Range("B10:E10").NumberFormat = "h:mm"

# Tuesday 2018-07-03 timesheet entry (partial)
$ws.Range("A11").Value = 43284
$ws.Range("B11").Value = 0.38194444444444442
$ws.Range("D11").Value = 0.3923611111111111
$ws.Range("B11").NumberFormat = "h:mm"
$ws.Range("D11").NumberFormat = "h:mm"

# Update active cell selection to G14
$ws.Range("G14").Select()
